$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) values are plain text in the source sheet (t="inlineStr").
# Writing numeric-looking strings straight to .Value lets Excel "helpfully"
# reinterpret them as numbers (stripping e.g. a trailing zero: "0.990" -> 0.99).
# Using .Formula with a leading text-qualifier apostrophe forces literal text,
# then ClearFormats() drops the transient @ (Text) number-format Excel applies
# so the cell keeps the workbooks original (unstyled) look.

# Row 2
$ws.Range("D2").Formula = "'57.274.75"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.41%  "

# Row 3
$ws.Range("D3").Formula = "'2.332.79"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.40%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").Formula = "'534.04"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.93%  "

# Row 6
$ws.Range("D6").Formula = "'133.07"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.05%  "

# Row 7
$ws.Range("E7").Value = "  -0.49%  "

# Row 8
$ws.Range("E8").Value = "  -0.50%  "

# Row 9
$ws.Range("D9").Formula = "'2.365.16"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.50%  "

# Row 10
$ws.Range("E10").Value = "  -0.95%  "

# Row 11
$ws.Range("E11").Value = "  +0.47%  "

# Row 12
$ws.Range("E12").Value = "  -1.61%  "

# Row 13
$ws.Range("D13").Formula = "'0.345"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.43%  "

# Row 14
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Formula = "'2.753.36"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.24%  "

# Row 15
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").Formula = "'23.62"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.43%  "

# Row 16
$ws.Range("D16").Formula = "'57.226.11"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.38%  "

# Row 17
$ws.Range("E17").Value = "  -1.10%  "

# Row 18
$ws.Range("D18").Formula = "'2.354.31"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.49%  "

# Row 19
$ws.Range("D19").Formula = "'339.92"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.97%  "

# Row 20
$ws.Range("D20").Formula = "'10.47"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.44%  "

# Row 21
$ws.Range("D21").Formula = "'6.93"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.62%  "

# Row 22
$ws.Range("E22").Value = "  -1.51%  "

# Row 23
$ws.Range("E23").Value = "  +0.14%  "

# Row 24
$ws.Range("D24").Formula = "'61.88"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.00%  "

# Row 25
$ws.Range("E25").Value = "  +10.01%  "

# Row 26
$ws.Range("E26").Value = "  -0.02%  "

# Row 27
$ws.Range("D27").Formula = "'0.990"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.69%  "

# Row 28
$ws.Range("E28").Value = "  +3.37%  "

# Row 29
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Formula = "'1.73"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.96%  "

# Row 30
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").Formula = "'169.89"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.11%  "

# Row 31
$ws.Range("E31").Value = "  -1.74%  "

# Row 32
$ws.Range("D32").Formula = "'6.14"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.77%  "

# Row 33
$ws.Range("E33").Value = "  -0.20%  "

# Row 34
$ws.Range("D34").Formula = "'0.998"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.06%  "

# Row 35
$ws.Range("E35").Value = "  -0.23%  "

# Row 36
$ws.Range("E36").Value = "  -0.61%  "

# Row 37
$ws.Range("E37").Value = "  +0.08%  "

# Row 38
$ws.Range("D38").Formula = "'0.912"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.31%  "

# Row 39
$ws.Range("E39").Value = "  +1.06%  "

# Row 40
$ws.Range("D40").Formula = "'38.98"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.31%  "

# Row 41
$ws.Range("D41").Formula = "'148.09"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.90%  "

# Row 42
$ws.Range("E42").Value = "  -1.61%  "

# Row 43
$ws.Range("D43").Formula = "'287.10"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.34%  "

# Row 44
$ws.Range("E44").Value = "  -1.61%  "

# Row 45
$ws.Range("E45").Value = "  -1.25%  "

# Row 46
$ws.Range("E46").Value = "  -0.33%  "

# Row 47
$ws.Range("D47").Formula = "'0.0505"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.36%  "

# Row 48
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Formula = "'0.563"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.03%  "

# Row 49
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Formula = "'18.82"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +4.54%  "

# Row 50
$ws.Range("E50").Value = "  -0.78%  "

# Row 51
$ws.Range("D51").Formula = "'17.39"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.01%  "
